$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC1, step 2 "Expected Results": append a sentence about the additional
# ordering of the list by the arrival date of the request at the
# authorization stage.
$ws.Range("D10").Value = "SYSTEM Recupera e exibe para o usuário a lista de diárias aptas para pagamento ordenado pelo numero de diarias em ordem crescente. Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de autorização (após registrar o empenho)."

# The "filter by user" test step/result (previously under TC4) now becomes
# the content of TC2's second step/result.
$ws.Range("B20").Value = "Chefe Seleciona um usuário para filtrar as autorizações de pagamento associadas a ele; e Submete a busca ao sistema."
$ws.Range("D20").Value = "SYSTEM Filtra os registros (autorizações de pagamento pendentes) e exibe apenas aqueles atribuídos ao usuário selecionado."

# The "assign/unassign responsibility" test step/result (previously TC2's
# content) shifts down to become TC3's second step/result.
$ws.Range("B28").Value = "Chefe Dado um registro selecionado (solicitação aguardando autorização de pagamento - AP), o usuário pode atribuir/desatribuir a responsabilidade da AP a si próprio; e Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D28").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela AP) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

# The "register payment authorization" test step/result (previously TC3's
# content) shifts down to become TC4's second step/result.
$ws.Range("B36").Value = "Chefe Clica para realizar a autorização de pagamento."
$ws.Range("D36").Value = "SYSTEM Apresenta a tela de Registrar Autorizações de Pagamento"
